$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update tStart (column B) values - all become 1E-3 (0.001)
$ws.Range("B2:B16").Value = 0.001

# Update tStop (column C) values per-row
$ws.Range("C2").Value = 0.0035
$ws.Range("C3").Value = 0.0035
$ws.Range("C4").Value = 0.0035
$ws.Range("C5").Value = 0.0035
$ws.Range("C6").Value = 0.0035
$ws.Range("C7").Value = 0.0022
$ws.Range("C8").Value = 0.0022
$ws.Range("C9").Value = 0.0022
$ws.Range("C10").Value = 0.0022
$ws.Range("C11").Value = 0.0035
$ws.Range("C12").Value = 0.0035
$ws.Range("C13").Value = 0.002
$ws.Range("C14").Value = 0.002
$ws.Range("C15").Value = 0.0035
$ws.Range("C16").Value = 0.0035

# Update the active selection to F7
$ws.Range("F7").Select()
